# postMessage with ExcelMessage as input
#
# Changes applied to the first sheet ("第一棵分類樹"):
#   - I1 header: "5_Url"   -> "6_Url"
#   - J1 header: "6_RefId" -> "5_RefId"
#   - D4:        "爭界"     -> "爭界123"
#   - Active selection moves from E1 to D4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("第一棵分類樹")

$ws.Range("J1").Value = "5_RefId"
$ws.Range("I1").Value = "6_Url"
$ws.Range("D4").Value = "爭界123"

$ws.Activate()
$ws.Range("D4").Select()
